$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new rows right after the header (at row 2), pushing existing
#    data down, and populate them with the new account records. The account
#    numbers are zero-padded digit strings, so force text formatting before
#    assigning the value (otherwise Excel auto-coerces them to numbers and
#    drops the leading zeros); ClearFormats() afterwards removes the
#    temporary "@" number-format style so the cell ends up unstyled, same as
#    every other "Conta" cell in the sheet.
$ws.Range("A2:A3").EntireRow.Insert()

$ws.Cells.Item(2,1).NumberFormat = "@"
$ws.Cells.Item(2,1).Value = "005646524"
$ws.Cells.Item(2,1).ClearFormats()
$ws.Cells.Item(2,2).Value = "EVANGELINA"
$ws.Cells.Item(2,3).Value = 53000

$ws.Cells.Item(3,1).NumberFormat = "@"
$ws.Cells.Item(3,1).Value = "000806386"
$ws.Cells.Item(3,1).ClearFormats()
$ws.Cells.Item(3,2).Value = "FERNANDA"
$ws.Cells.Item(3,3).Value = 36457.46

# 2. The row that used to hold "005341184 / BRENO" (originally row 9) is now
#    at row 11 after the two-row insert above; rename the account holder
#    while keeping the balance unchanged.
$ws.Cells.Item(11,1).NumberFormat = "@"
$ws.Cells.Item(11,1).Value = "005295509"
$ws.Cells.Item(11,1).ClearFormats()
$ws.Cells.Item(11,2).Value = "BHRUNA"

# 3. The old duplicate "000806386 / FERNANDA / 457.46" row (originally row 20)
#    is now at row 22; remove it entirely since FERNANDA's correct balance is
#    already captured in the newly inserted row above.
$ws.Rows(22).Delete()
